# Apply the "updated rmi files 3.4.3" edit:
#   - About sheet: add region label "Minnesota" (B1) and a last-updated date (C1)
#   - GbPbT sheet: hard-code the region-specific pollutant GWP values (VOC, CO, NOx,
#     PM10, PM25, SOx, BC, OC) to 0 instead of pulling them from the Data tab,
#     and make GbPbT the active/selected sheet.

$wb = $excel.ActiveWorkbook

# --- About sheet -----------------------------------------------------------
$about = $wb.Worksheets.Item("About")

$about.Range("B1").Value = "Minnesota"

$about.Range("C1").Value = 44861
$about.Range("C1").NumberFormat = "mm-dd-yy"

# --- GbPbT sheet -------------------------------------------------------------
$gbpbt = $wb.Worksheets.Item("GbPbT")

$rows = 3, 4, 5, 6, 7, 8, 9, 10
foreach ($r in $rows) {
    $gbpbt.Cells.Item($r, 2).Value = 0
    $gbpbt.Cells.Item($r, 3).Value = 0
}

# GbPbT becomes the selected/active sheet (tab focus moved from About to GbPbT).
$gbpbt.Activate()
$gbpbt.Range("F9").Select()
